$d = $word.ActiveDocument

# The table for the "4η ΕΒΔΟΜΑΔΑ" (17/10 - 23/10) weekly menu grid is the
# 4th table in the document. Its 6th column (20/10, currently 1683 twips /
# 84.15 pt) and 7th column (21/10, currently 1706 twips / 85.3 pt) are being
# resized so the boundary between them shifts right: the 6th column grows to
# 2209 twips (110.45 pt) and the 7th column shrinks to 1180 twips (59 pt).
# The combined width of the two columns (3389 twips) is unchanged.
$t = $d.Tables.Item(4)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $t.Cell($r, 6).Width = 110.45
    $t.Cell($r, 7).Width = 59.0
}
